$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.959.32'
$ws.Cells.Item(2, 5).Value = '  -6.21%  '

$ws.Cells.Item(3, 4).Value = '3.072.36'
$ws.Cells.Item(3, 5).Value = '  -9.02%  '

$ws.Cells.Item(4, 5).Value = '  -0.18%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '501.93'
$ws.Cells.Item(5, 5).Value = '  -4.65%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '164.16'
$ws.Cells.Item(6, 5).Value = '  -11.64%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.573'
$ws.Cells.Item(7, 5).Value = '  -4.45%  '

$ws.Cells.Item(8, 5).Value = '  -0.12%  '

$ws.Cells.Item(9, 4).Value = '3.073.10'
$ws.Cells.Item(9, 5).Value = '  -8.66%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.572'
$ws.Cells.Item(10, 5).Value = '  -8.14%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '50.93'
$ws.Cells.Item(11, 5).Value = '  -12.18%  '

$ws.Cells.Item(12, 5).Value = '  -6.30%  '

$ws.Cells.Item(13, 5).Value = '  -5.35%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '8.55'
$ws.Cells.Item(14, 5).Value = '  -7.10%  '

$ws.Cells.Item(15, 4).Value = '3.574.58'
$ws.Cells.Item(15, 5).Value = '  -8.75%  '

$ws.Cells.Item(16, 5).Value = '  -9.39%  '

$ws.Cells.Item(17, 4).Value = '3.081.42'
$ws.Cells.Item(17, 5).Value = '  -8.80%  '

$ws.Cells.Item(18, 4).Value = '60.904.93'
$ws.Cells.Item(18, 5).Value = '  -5.94%  '

$ws.Cells.Item(19, 5).Value = '  -5.47%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '10.39'
$ws.Cells.Item(20, 5).Value = '  -6.32%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.915'
$ws.Cells.Item(21, 5).Value = '  -5.41%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '350.38'
$ws.Cells.Item(22, 5).Value = '  -5.42%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '77.52'
$ws.Cells.Item(23, 5).Value = '  -4.32%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.54'
$ws.Cells.Item(24, 5).Value = '  -4.51%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '6.04'
$ws.Cells.Item(25, 5).Value = '  +4.00%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '10.50'
$ws.Cells.Item(26, 5).Value = '  -2.71%  '

$ws.Cells.Item(27, 5).Value = '  -0.40%  '

$ws.Cells.Item(28, 5).Value = '  -5.39%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '10.63'
$ws.Cells.Item(29, 5).Value = '  -6.85%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.72'
$ws.Cells.Item(30, 5).Value = '  -9.11%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '620.30'
$ws.Cells.Item(31, 5).Value = '  -6.79%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '26.95'
$ws.Cells.Item(32, 5).Value = '  -8.50%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.06'
$ws.Cells.Item(33, 5).Value = '  -9.36%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '10.74'
$ws.Cells.Item(34, 5).Value = '  -3.67%  '

$ws.Cells.Item(35, 2).Value = 'Dai'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.02%  '

$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0989'
$ws.Cells.Item(36, 5).Value = '  -5.96%  '

$ws.Cells.Item(37, 2).Value = 'OKB'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '54.48'
$ws.Cells.Item(37, 5).Value = '  -10.82%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '34.55'
$ws.Cells.Item(38, 5).Value = '  -5.08%  '

$ws.Cells.Item(39, 2).Value = 'TheGraph'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.357'
$ws.Cells.Item(39, 5).Value = '  -5.36%  '

$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.999'
$ws.Cells.Item(40, 5).Value = '  -0.08%  '

$ws.Cells.Item(41, 4).Value = '0.0₃0648'
$ws.Cells.Item(41, 5).Value = '  +4.56%  '

$ws.Cells.Item(42, 5).Value = '  -7.32%  '

$ws.Cells.Item(43, 4).Value = '2.743.23'
$ws.Cells.Item(43, 5).Value = '  -3.59%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.37'
$ws.Cells.Item(44, 5).Value = '  +2.04%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.56'
$ws.Cells.Item(45, 5).Value = '  -3.49%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.78'
$ws.Cells.Item(46, 5).Value = '  +8.29%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0369'
$ws.Cells.Item(47, 5).Value = '  -5.61%  '

$ws.Cells.Item(48, 5).Value = '  -11.30%  '

$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.118'
$ws.Cells.Item(49, 5).Value = '  -4.88%  '

$ws.Cells.Item(50, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.80'
$ws.Cells.Item(50, 5).Value = '  -1.47%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '128.28'
$ws.Cells.Item(51, 5).Value = '  -6.70%  '
